# Generate Report for Handoff
# Update the "Latest Handoff Date"/"Latest Handoff Datetime" values recorded
# for the a81892bb-1ba4-45a6-a482-626111288396 file to reflect the new
# handoff run timestamp.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 7 is the a81892bb-... file; column D = "Latest Handoff Date"
$overview.Range("D7").Value = "2016-03-13 21:03:08"

# zh-cn sheet: row 7 is the a81892bb-... file; column E = "Latest Handoff Datetime"
$zhcn.Range("E7").Value = "2016-03-13 21:03:05"

# de-de sheet: row 7 is the a81892bb-... file; column E = "Latest Handoff Datetime"
$dede.Range("E7").Value = "2016-03-13 21:03:08"
